# Generate Report for Handoff
# Swap the d2f80547 / d4bd5148 rows (the "Overview" row 3 moves to "Ready for
# handoff" with a new handoff time, and the per-language rows that used to
# describe d2f80547 now describe d4bd5148 and vice versa), then refresh the
# hyperlinks so their (unchanged) targets keep the same ordering while their
# display text follows the new row contents.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-23 09:39:23"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
$ws2.Range("D2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf"
$ws2.Range("F2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
$ws2.Range("G2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf"

$ws2.Range("A3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-23 09:39:16"
$ws2.Range("F3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
$ws2.Range("G3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad925be95dc0b9e0e12e76e6756ebc9e085395cf/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/5b75865346e440d40aef8a0c26733db409561046/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/495b27d90da480cde4011d7bb5013f04d6cc7369/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad925be95dc0b9e0e12e76e6756ebc9e085395cf/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/5b75865346e440d40aef8a0c26733db409561046/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/495b27d90da480cde4011d7bb5013f04d6cc7369/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
$ws3.Range("D2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf"
$ws3.Range("F2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
$ws3.Range("G2").Value = "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf"

$ws3.Range("A3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-23 09:39:23"
$ws3.Range("F3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
$ws3.Range("G3").Value = "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8960772060e52c3708254bb587723ac2fb69fd6b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/2005b4baec89f7cd0051cd3393e7ee21da81365e/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/96d8feafe5230bc6c6484ace31302012a5eba910/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf", "", "", "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8960772060e52c3708254bb587723ac2fb69fd6b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/2005b4baec89f7cd0051cd3393e7ee21da81365e/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/96d8feafe5230bc6c6484ace31302012a5eba910/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf") | Out-Null

$wb.Save()
